# Natmi following Dr Hou advice
# Rebuild the LR-pairs result table (Vip -> Vipr1) with the updated
# sending/target cluster combinations and recomputed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out any previously-written data rows below the header so the
# sheet only contains the header row plus the new result rows.
$ws.Range("A2:T11").ClearContents()

$data = New-Object "object[,]" 10,20
# Row 2
$data[0,0] = "ECs"
$data[0,1] = "Vip"
$data[0,2] = "Vipr1"
$data[0,3] = "FAPs"
$data[0,4] = 2
$data[0,5] = 1
$data[0,6] = 1.38698
$data[0,7] = 2.77396
$data[0,8] = 0.5967131687073423
$data[0,9] = 0.4965808999056411
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.032709
$data[0,13] = 0.098127
$data[0,14] = 0.006659324197390494
$data[0,15] = 0.006914300080637592
$data[0,16] = 0.04536672882
$data[0,17] = 0.27220037292
$data[0,18] = 0.00397370644327436
$data[0,19] = 0.003433509356260662
# Row 3
$data[1,0] = "ECs"
$data[1,1] = "Vip"
$data[1,2] = "Vipr1"
$data[1,3] = "M1"
$data[1,4] = 2
$data[1,5] = 1
$data[1,6] = 1.38698
$data[1,7] = 2.77396
$data[1,8] = 0.5967131687073423
$data[1,9] = 0.4965808999056411
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 1.571286666666667
$data[1,13] = 4.71386
$data[1,14] = 0.3199030028545778
$data[1,15] = 0.3321516257310865
$data[1,16] = 2.179343180933333
$data[1,17] = 13.0760590856
$data[1,18] = 0.1908903345123491
$data[1,19] = 0.1649401532106646
# Row 4
$data[2,0] = "ECs"
$data[2,1] = "Vip"
$data[2,2] = "Vipr1"
$data[2,3] = "M2"
$data[2,4] = 2
$data[2,5] = 1
$data[2,6] = 1.38698
$data[2,7] = 2.77396
$data[2,8] = 0.5967131687073423
$data[2,9] = 0.4965808999056411
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 2.432594333333333
$data[2,13] = 7.297783
$data[2,14] = 0.4952592346571789
$data[2,15] = 0.5142219937975853
$data[2,16] = 3.373959688446666
$data[2,17] = 20.24375813068
$data[2,18] = 0.2955277072438584
$data[2,19] = 0.2553528204312779
# Row 5
$data[3,0] = "ECs"
$data[3,1] = "Vip"
$data[3,2] = "Vipr1"
$data[3,3] = "Neutro"
$data[3,4] = 2
$data[3,5] = 1
$data[3,6] = 1.38698
$data[3,7] = 2.77396
$data[3,8] = 0.5967131687073423
$data[3,9] = 0.4965808999056411
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.3317826666666666
$data[3,13] = 0.9953479999999999
$data[3,14] = 0.06754863616766264
$data[3,15] = 0.07013497566075051
$data[3,16] = 0.4601759230133332
$data[3,17] = 2.761055538079999
$data[3,18] = 0.04030716072946536
$data[3,19] = 0.03482768932847572
# Row 6
$data[4,0] = "ECs"
$data[4,1] = "Vip"
$data[4,2] = "Vipr1"
$data[4,3] = "sCs"
$data[4,4] = 2
$data[4,5] = 1
$data[4,6] = 1.38698
$data[4,7] = 2.77396
$data[4,8] = 0.5967131687073423
$data[4,9] = 0.4965808999056411
$data[4,10] = 2
$data[4,11] = 1
$data[4,12] = 0.5433870000000001
$data[4,13] = 1.086774
$data[4,14] = 0.1106298021231902
$data[4,15] = 0.07657710472994017
$data[4,16] = 0.75366690126
$data[4,17] = 3.01466760504
$data[4,18] = 0.06601425977839509
$data[4,19] = 0.03802672757896221
# Row 7
$data[5,0] = "Neutro"
$data[5,1] = "Vip"
$data[5,2] = "Vipr1"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.9373863333333334
$data[5,7] = 2.812159
$data[5,8] = 0.4032868312926577
$data[5,9] = 0.503419100094359
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.032709
$data[5,13] = 0.098127
$data[5,14] = 0.006659324197390494
$data[5,15] = 0.006914300080637592
$data[5,16] = 0.03066096957700001
$data[5,17] = 0.275948726193
$data[5,18] = 0.002685617754116134
$data[5,19] = 0.00348079072437693
# Row 8
$data[6,0] = "Neutro"
$data[6,1] = "Vip"
$data[6,2] = "Vipr1"
$data[6,3] = "M1"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.9373863333333334
$data[6,7] = 2.812159
$data[6,8] = 0.4032868312926577
$data[6,9] = 0.503419100094359
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 1.571286666666667
$data[6,13] = 4.71386
$data[6,14] = 0.3199030028545778
$data[6,15] = 0.3321516257310865
$data[6,16] = 1.472902647082222
$data[6,17] = 13.25612382374
$data[6,18] = 0.1290126683422287
$data[6,19] = 0.1672114725204219
# Row 9
$data[7,0] = "Neutro"
$data[7,1] = "Vip"
$data[7,2] = "Vipr1"
$data[7,3] = "M2"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.9373863333333334
$data[7,7] = 2.812159
$data[7,8] = 0.4032868312926577
$data[7,9] = 0.503419100094359
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.432594333333333
$data[7,13] = 7.297783
$data[7,14] = 0.4952592346571789
$data[7,15] = 0.5142219937975853
$data[7,16] = 2.280280682610778
$data[7,17] = 20.522526143497
$data[7,18] = 0.1997315274133205
$data[7,19] = 0.2588691733663074
# Row 10
$data[8,0] = "Neutro"
$data[8,1] = "Vip"
$data[8,2] = "Vipr1"
$data[8,3] = "Neutro"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.9373863333333334
$data[8,7] = 2.812159
$data[8,8] = 0.4032868312926577
$data[8,9] = 0.503419100094359
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.3317826666666666
$data[8,13] = 0.9953479999999999
$data[8,14] = 0.06754863616766264
$data[8,15] = 0.07013497566075051
$data[8,16] = 0.3110085373702222
$data[8,17] = 2.799076836332
$data[8,18] = 0.02724147543819728
$data[8,19] = 0.03530728633227479
# Row 11
$data[9,0] = "Neutro"
$data[9,1] = "Vip"
$data[9,2] = "Vipr1"
$data[9,3] = "sCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.9373863333333334
$data[9,7] = 2.812159
$data[9,8] = 0.4032868312926577
$data[9,9] = 0.503419100094359
$data[9,10] = 2
$data[9,11] = 1
$data[9,12] = 0.5433870000000001
$data[9,13] = 1.086774
$data[9,14] = 0.1106298021231902
$data[9,15] = 0.07657710472994017
$data[9,16] = 0.5093635475110001
$data[9,17] = 3.056181285066001
$data[9,18] = 0.04461554234479512
$data[9,19] = 0.03855037715097796

$ws.Range("A2:T11").Value = $data

Write-Output "Updated Vip-Vipr1 LR-pairs rows 2..11 (A:T)"
